$wb = $excel.ActiveWorkbook

# Add the new "stage" worksheet at the end of the workbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "stage"

# Move it after the last existing sheet ("schedule")
$scheduleSheet = $wb.Worksheets.Item("schedule")
$newSheet.Move($null, $scheduleSheet)

$ws = $wb.Worksheets.Item("stage")

# Header row (row 2)
$ws.Range("A2").Value = "title"
$ws.Range("B2").Value = "backgroundColor"
$ws.Range("C2").Value = "dueDate"
$ws.Range("D2").Value = "icon"
$ws.Range("E2").Value = "internalComment"
$ws.Range("F2").Value = "isCrossLane"
$ws.Range("G2").Value = "isDeleted"
$ws.Range("H2").Value = "timetableOrder"
$ws.Range("I2").Value = "isReadOnly"
$ws.Range("J2").Value = "textColor"
$ws.Range("K2").Value = "type"
$ws.Range("L2").Value = "viewOrder"

# Example data row (row 3)
$ws.Range("A3").Value = "Main Stage"
$ws.Range("E3").Value = "Beispiel Main Stage"
$ws.Range("F3").Value = "false"
$ws.Range("G3").Value = "false"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = "false"
$ws.Range("L3").Value = 1
